# Update "想去人数" (want-to-go count) figures for two exhibitions
# on both the "展览" sheet and the "全部类型" sheet, which mirror the
# same underlying data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 346
    $ws.Range("F6").Value = 47
    $ws.Range("F9").Value = 343
}
